$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 64
$ws.Range("H64").Value = 4195.4165
$ws.Range("I64").Value = 3620.7144
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3620.7144
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -3372.7144
$ws.Range("N64").Value = -5496

# Row 67
$ws.Range("H67").Value = 4195.4165
$ws.Range("I67").Value = 3620.7144
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3620.7144
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2762.7144
$ws.Range("N67").Value = -6716

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 7794.3203
$ws.Range("I32").Value = 5909.3384
$ws.Range("K32").Value = 5909.3384
$ws.Range("M32").Value = -5622.3384

# Row 45
$ws.Range("H45").Value = 1694.2222
$ws.Range("I45").Value = 1435.4286
$ws.Range("J45").Value = 2600
$ws.Range("K45").Value = 1435.4286
$ws.Range("L45").Value = 2600
$ws.Range("M45").Value = -1058.4286
$ws.Range("N45").Value = -3354

# Row 61
$ws.Range("H61").Value = 3952.5193
$ws.Range("I61").Value = 2873.7437
$ws.Range("J61").Value = 7188.846
$ws.Range("K61").Value = 2873.7437
$ws.Range("L61").Value = 7188.846
$ws.Range("M61").Value = -2661.7437
$ws.Range("N61").Value = -7612.846

# Row 63 - M/N cells are removed entirely (row total columns collapse to L)
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()

# Row 66 - M/N cells are removed entirely (row total columns collapse to L)
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

# Row 122
$ws.Range("H122").Value = 15628758
$ws.Range("I122").Value = 3690.5
$ws.Range("K122").Value = 11071.5
$ws.Range("M122").Value = -8621.5

# Row 132
$ws.Range("H132").Value = 6299.9033
$ws.Range("I132").Value = 1521.2858
$ws.Range("J132").Value = 10235.235
$ws.Range("K132").Value = 4563.857400000001
$ws.Range("L132").Value = 30705.705
$ws.Range("M132").Value = -2033.857400000001
$ws.Range("N132").Value = -35765.705

# Row 136
$ws.Range("H136").Value = 3952.5193
$ws.Range("I136").Value = 2873.7437
$ws.Range("J136").Value = 7188.846
$ws.Range("K136").Value = 8621.231100000001
$ws.Range("L136").Value = 21566.538
$ws.Range("M136").Value = -6071.231100000001
$ws.Range("N136").Value = -26666.538

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 105
$ws.Range("H105").Value = 5110.311
$ws.Range("I105").Value = 4745.1
$ws.Range("J105").Value = 5840.7334
$ws.Range("K105").Value = 4745.1
$ws.Range("L105").Value = 5840.7334
$ws.Range("M105").Value = -2998.1
$ws.Range("N105").Value = -9334.733400000001

# Row 134
$ws.Range("H134").Value = 2237.3
$ws.Range("I134").Value = 2115.1177
$ws.Range("J134").Value = 2929.6667
$ws.Range("K134").Value = 6345.353099999999
$ws.Range("L134").Value = 8789.000100000001
$ws.Range("M134").Value = -3810.353099999999
$ws.Range("N134").Value = -13859.0001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 4
$ws.Range("H4").Value = 8355.777
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 10171.714
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 10171.714
$ws.Range("M4").Value = -1888
$ws.Range("N4").Value = -10395.714

# Row 18 (no M18 cell, both before and after)
$ws.Range("H18").Value = 55870.4
$ws.Range("J18").Value = 55870.4
$ws.Range("L18").Value = 55870.4
$ws.Range("N18").Value = -56330.4

# Row 31
$ws.Range("H31").Value = 2569.7908
$ws.Range("I31").Value = 1870.7106
$ws.Range("J31").Value = 7882.8
$ws.Range("K31").Value = 1870.7106
$ws.Range("L31").Value = 7882.8
$ws.Range("M31").Value = -1575.7106
$ws.Range("N31").Value = -8472.799999999999

# Row 34
$ws.Range("H34").Value = 2569.7908
$ws.Range("I34").Value = 1870.7106
$ws.Range("J34").Value = 7882.8
$ws.Range("K34").Value = 1870.7106
$ws.Range("L34").Value = 7882.8
$ws.Range("M34").Value = -1668.7106
$ws.Range("N34").Value = -8286.799999999999

# Row 132
$ws.Range("H132").Value = 3427.7917
$ws.Range("I132").Value = 2862.5
$ws.Range("J132").Value = 5123.6665
$ws.Range("K132").Value = 8587.5
$ws.Range("L132").Value = 15370.9995
$ws.Range("M132").Value = -6057.5
$ws.Range("N132").Value = -20430.9995

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 4 (no J4 change, no L4/N4 here - not present in diff for this row's hunk)
$ws.Range("H4").Value = 581.5833
$ws.Range("I4").Value = 247.375
$ws.Range("K4").Value = 742.125
$ws.Range("M4").Value = -630.125

# Row 12
$ws.Range("H12").Value = 25641264
$ws.Range("J12").Value = 238.11539
$ws.Range("L12").Value = 714.34617
$ws.Range("N12").Value = -1060.34617

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 87 - N87 cell removed entirely
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90 - N90 cell removed entirely
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 122
$ws.Range("H122").Value = 16875
$ws.Range("I122").Value = 26250
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 78750
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -76300
$ws.Range("N122").Value = -27400

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 2
$ws.Range("H2").Value = 100000
$ws.Range("J2").Value = 100000
$ws.Range("L2").Value = 100000
$ws.Range("N2").Value = -100224

# Row 30
$ws.Range("H30").Value = 8338.666999999999
$ws.Range("I30").Value = 5016
$ws.Range("K30").Value = 5016
$ws.Range("M30").Value = -4908

# Row 122
$ws.Range("H122").Value = 5665.1333
$ws.Range("I122").Value = 4472.643
$ws.Range("J122").Value = 7629.2354
$ws.Range("K122").Value = 13417.929
$ws.Range("L122").Value = 22887.7062
$ws.Range("M122").Value = -10967.929
$ws.Range("N122").Value = -27787.7062

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 3608.6553
$ws.Range("I122").Value = 1689.2084
$ws.Range("J122").Value = 12822
$ws.Range("K122").Value = 5067.6252
$ws.Range("L122").Value = 38466
$ws.Range("M122").Value = -2617.6252
$ws.Range("N122").Value = -43366

# Row 132
$ws.Range("H132").Value = 3218.3125
$ws.Range("I132").Value = 2699.0833
$ws.Range("J132").Value = 4776
$ws.Range("K132").Value = 8097.249899999999
$ws.Range("L132").Value = 14328
$ws.Range("M132").Value = -5567.249899999999
$ws.Range("N132").Value = -19388
